# Dedication paragraph: "For Rosie, Joey and Zach, in the hope ..."
# becomes: "For Rosie, Joey, Zach, Amy, Elliott and Lyla, in the hope ..."
# with Amy / Elliott / Lyla bolded (matching Rosie / Joey / Zach).

$d = $word.ActiveDocument

# Locate the dedication paragraph robustly (rather than hard-coding an index).
$p = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -like "*Joey and Zach*") {
        $p = $cand
        break
    }
}
if ($p -eq $null) {
    throw "Could not find the dedication paragraph"
}

# --- Step 1: "Joey and Zach" -> "Joey, Zach" ------------------------------
# i.e. the run containing the space right after "Joey" becomes "," and the
# following "and" run is removed, leaving the existing space before "Zach".
$r = $p.Range.Duplicate
$r.Find.Execute(" and Zach", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$spaceRange = $d.Range($r.Start, $r.Start + 1)       # the space right after "Joey"
$andRange   = $d.Range($r.Start + 1, $r.Start + 4)   # "and"
$spaceRange.Text = ","
$andRange.Delete()

# --- Step 2: insert the plain-text skeleton after "Zach" ------------------
# Splitting inside the existing (non-bold) ", in the hope..." run keeps the
# inserted punctuation/spacing runs italic-only (no stray bold bleed-through).
$rZach = $p.Range.Duplicate
$rZach.Find.Execute("Zach", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$afterZach = $rZach.End
$insPos = $d.Range($afterZach + 2, $afterZach + 2)   # 2 chars in = right after ", "
$insPos.InsertBefore("Amy, Elliott and Lyla, ")

# --- Step 3: turn the plain "Amy" / "Elliott" / "Lyla" placeholders into ---
# bold+italic runs matching Rosie/Joey/Zach's formatting (including bCs,
# which the Bold property alone does not set in this runtime -- so an
# existing bold run is copy/pasted as a template and then renamed in place).
function Make-Bold-Name($paragraph, $name) {
    $rr = $paragraph.Range.Duplicate
    $rr.Find.Execute($name, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
    $start = $rr.Start
    $end = $rr.End

    $tmpl = $paragraph.Range.Duplicate
    $tmpl.Find.Execute("Zach", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
    $tmpl.Copy()

    $dest = $d.Range($start, $start)
    $dest.Paste()

    # the pasted template is "Zach" (4 chars); the old plain name text is now
    # shifted forward by 4 characters
    $oldStart = $start + 4
    $oldEnd = $end + 4
    $toDelete = $d.Range($oldStart, $oldEnd)
    $toDelete.Delete()

    $pasted = $d.Range($start, $start + 4)
    $pasted.Text = $name
}

Make-Bold-Name $p "Amy"
Make-Bold-Name $p "Elliott"
Make-Bold-Name $p "Lyla"
